# Rutuja Added test case for update profile
#
# Adds a new "TestData-Rutuja" worksheet (with a small name-splitting scratch
# table) after the existing three sheets, makes it the active tab, and nudges
# a couple of other sheets' selections/zoom back to their "closed" state.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last tab -------------------------------
$lastIndex = $wb.Worksheets.Count
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws4.Name = "TestData-Rutuja"

# --- Fill in the scratch data ---------------------------------------------
# Write order matters for shared-string allocation order (A8/Jadhav is
# written before B7/Rutuja so the workbook's shared-string table comes out
# in the same order as the authored file).
$ws4.Range("A1").Value = "utuja"
$ws4.Range("A6").Value = "r"
$ws4.Range("B6").Value = "Rj"
$ws4.Range("A7").Value = "Rutu"
$ws4.Range("A8").Value = "Jadhav"
$ws4.Range("B7").Value = "Rutuja"
$ws4.Range("B8").Value = "J"

# Leave the cursor on B8, as in the final authored state.
[void]$ws4.Range("B8").Select()

# --- Registration_Details: selection moved to A7 --------------------------
$ws3 = $wb.Worksheets.Item("Registration_Details")
[void]$ws3.Range("A7").Select()

# --- Credentials: zoom reset back to 100% ----------------------------------
$ws1 = $wb.Worksheets.Item("Credentials")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 100

# --- Make the new sheet the active tab (last-activated wins) --------------
$ws4.Activate()
